{"js": "// Find the paragraph that ends with the \"Using historical data...\" sentence\n// about McDonald's stock forecasts, then insert a new \"source code\" style\n// paragraph right after it (and before the forecasts table), containing:\n//   pander(table_forecasts)\n// split into two runs: \"pander\" (FunctionTok char style) and\n// \"(table_forecasts)\" (NormalTok char style).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"forecasts are summarized in the table below\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph containing marker text.\");\n}\n\n// Insert a new empty paragraph right after the target paragraph.\nconst newPara = target.insertParagraph(\"\", \"After\");\nnewPara.style = \"Source Code\";\n\n// First run: \"pander\" with FunctionTok character style.\nconst r1 = newPara.insertText(\"pander\", \"End\");\nr1.style = \"FunctionTok\";\n\n// Second run: \"(table_forecasts)\" with NormalTok character style.\nconst r2 = newPara.insertText(\"(table_forecasts)\", \"End\");\nr2.style = \"NormalTok\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"Source Code\" style paragraph, containing the rendered\n# pander() call that produces the forecasts table, right after the\n# \"Using historical data...\" paragraph and right before the table itself.\n\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"forecasts are summarized in the table below\")\nif (-not $found) {\n    throw \"Could not find target paragraph containing marker text.\"\n}\n\n$targetPara = $findRange.Paragraphs(1)\n$targetPara.Range.InsertParagraphAfter()\n\n$newPara = $targetPara.Next()\n$newRange = $newPara.Range\n$newRange.Style = \"Source Code\"\n\n$insertPoint = $newRange.Start\n\n# First run: \"pander\" with FunctionTok character style.\n$run1 = $d.Range($insertPoint, $insertPoint)\n$run1.InsertAfter(\"pander\")\n$run1.Style = \"FunctionTok\"\n\n# Second run: \"(table_forecasts)\" with NormalTok character style.\n$run2 = $d.Range($run1.End, $run1.End)\n$run2.InsertAfter(\"(table_forecasts)\")\n$run2.Style = \"NormalTok\"\n"}
